$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.484.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.32%  '
$ws.Range("D3").Value = "'3.408.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.31%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = "'579.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.90%  '
$ws.Range("D6").Value = "'133.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.46%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = "'3.403.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.48%  '
$ws.Range("D9").Value = "'0.479"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.58%  '
$ws.Range("D10").Value = "'0.119"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.64%  '
$ws.Range("D11").Value = "'7.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.06%  '
$ws.Range("D12").Value = "'0.371"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.04%  '
$ws.Range("D13").Value = "'3.991.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("D14").Value = "'0.0000177"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.66%  '
$ws.Range("D15").Value = "'3.437.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("D17").Value = "'25.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.74%  '
$ws.Range("D18").Value = "'64.558.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.00%  '
$ws.Range("D19").Value = "'9.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -14.90%  '
$ws.Range("D20").Value = "'5.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.60%  '
$ws.Range("D21").Value = "'13.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.46%  '
$ws.Range("D22").Value = "'377.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.16%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = "'0.536"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.74%  '
$ws.Range("D25").Value = "'71.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.14%  '
$ws.Range("D26").Value = "'3.551.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.09%  '
$ws.Range("D27").Value = "'0.0000102"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.04%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").Value = "'7.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.34%  '
$ws.Range("D30").Value = "'2.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.85%  '
$ws.Range("D31").Value = "'7.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.11%  '
$ws.Range("D32").Value = "'3.430.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.82%  '
$ws.Range("D36").Value = "'168.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.89%  '
$ws.Range("E37").Value = '  -13.49%  '
$ws.Range("D38").Value = "'6.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -12.50%  '
$ws.Range("D39").Value = "'1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -12.53%  '
$ws.Range("D40").Value = "'4.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -13.78%  '
$ws.Range("D41").Value = "'0.0752"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.98%  '
$ws.Range("D42").Value = "'0.799"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.02%  '
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("D44").Value = "'41.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.00%  '
$ws.Range("D45").Value = "'4.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -14.91%  '
$ws.Range("D46").Value = "'1.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.43%  '
$ws.Range("D47").Value = "'1.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("D48").Value = "'22.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.66%  '
$ws.Range("D49").Value = "'6.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.07%  '
$ws.Range("D50").Value = "'2.178.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.14%  '
$ws.Range("D51").Value = "'1.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -19.02%  '

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'22.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.84%  '
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.140"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.14%  '
